$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to
# Text format first, so Excel keeps them as strings (matching the source
# data, which stores all Price/Volume columns as text).
$textCells = @("D5", "D6", "D14", "D19", "D20", "D21", "D22", "D23", "D25", "D26", "D27", "D32", "D35", "D38", "D39", "D40", "D41", "D42", "D43", "D45", "D46", "D48", "D50")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = "62.456.76"
$ws.Range("E2").Value = "  -2.15%  "
$ws.Range("D3").Value = "2.445.23"
$ws.Range("E3").Value = "  -1.32%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "582.77"
$ws.Range("E5").Value = "  +1.03%  "
$ws.Range("D6").Value = "143.46"
$ws.Range("E6").Value = "  -3.90%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  -2.82%  "
$ws.Range("D9").Value = "2.441.72"
$ws.Range("E9").Value = "  -1.39%  "
$ws.Range("E10").Value = "  -4.84%  "
$ws.Range("E12").Value = "  -1.45%  "
$ws.Range("E13").Value = "  -3.86%  "
$ws.Range("D14").Value = "26.47"
$ws.Range("E14").Value = "  -3.40%  "
$ws.Range("E15").Value = "  -5.34%  "
$ws.Range("D16").Value = "2.859.74"
$ws.Range("E16").Value = "  -2.26%  "
$ws.Range("D17").Value = "62.279.25"
$ws.Range("E17").Value = "  -2.17%  "
$ws.Range("D18").Value = "2.436.45"
$ws.Range("E18").Value = "  -1.33%  "
$ws.Range("D19").Value = "11.01"
$ws.Range("E19").Value = "  -4.83%  "
$ws.Range("D20").Value = "7.13"
$ws.Range("E20").Value = "  -4.00%  "
$ws.Range("D21").Value = "332.38"
$ws.Range("E21").Value = "  +0.35%  "
$ws.Range("D22").Value = "4.12"
$ws.Range("E22").Value = "  -2.61%  "
$ws.Range("D23").Value = "1.96"
$ws.Range("E23").Value = "  -8.03%  "
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("D25").Value = "65.92"
$ws.Range("E25").Value = "  -0.53%  "
$ws.Range("D26").Value = "638.87"
$ws.Range("E26").Value = "  +0.72%  "
$ws.Range("D27").Value = "9.22"
$ws.Range("E27").Value = "  +4.32%  "
$ws.Range("D28").Value = "0.0₃0960"
$ws.Range("E28").Value = "  -9.52%  "
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("E31").Value = "  -6.85%  "
$ws.Range("D32").Value = "8.07"
$ws.Range("E32").Value = "  -3.92%  "
$ws.Range("E33").Value = "  -2.80%  "
$ws.Range("E34").Value = "  -1.39%  "
$ws.Range("D35").Value = "4.98"
$ws.Range("E35").Value = "  -6.13%  "
$ws.Range("E36").Value = "  +0.17%  "
$ws.Range("E37").Value = "  -5.82%  "
$ws.Range("D38").Value = "0.376"
$ws.Range("E38").Value = "  -2.28%  "
$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D39").Value = "149.79"
$ws.Range("E39").Value = "  +1.61%  "
$ws.Range("B40").Value = "EthereumClassic"
$ws.Range("C40").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D40").Value = "18.43"
$ws.Range("E40").Value = "  -2.47%  "
$ws.Range("D41").Value = "5.27"
$ws.Range("E41").Value = "  -4.56%  "
$ws.Range("D42").Value = "1.75"
$ws.Range("E42").Value = "  -4.20%  "
$ws.Range("D43").Value = "42.75"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").Value = "2.50"
$ws.Range("E45").Value = "  -9.65%  "
$ws.Range("D46").Value = "143.79"
$ws.Range("E46").Value = "  -5.10%  "
$ws.Range("E47").Value = "  -3.36%  "
$ws.Range("D48").Value = "0.0523"
$ws.Range("E48").Value = "  -4.03%  "
$ws.Range("E49").Value = "  -1.51%  "
$ws.Range("D50").Value = "19.73"
$ws.Range("E50").Value = "  -8.58%  "
$ws.Range("D51").Value = "0.0₆0232"
$ws.Range("E51").Value = "  +1.39%  "

# Restore default (unstyled) cell style so formatting matches the original
# workbook, which left these data cells without an explicit style index.
foreach ($c in $textCells) {
    $ws.Range($c).Style = "Normal"
}
